$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $expectedOld, $newText) {
    $cell = $t.Rows.Item($row).Cells.Item($col)
    $r = $cell.Range
    # Cell.Range.Text carries trailing end-of-cell / paragraph markers
    # (\r, \a) that aren't part of the visible text - trim them before
    # sanity-checking against the expected old value.
    $cur = $r.Text.TrimEnd([char]13, [char]7)
    if ($cur -ne $expectedOld) {
        Write-Host "WARNING: row $row col $col expected '$expectedOld' but found '$cur'"
    }
    $r.Text = $newText
}

# CP row (N)
Set-CellText 116 2 "218,853" "218,854"

# Censored row
Set-CellText 117 3 "207,378 (95%)" "207,083 (95%)"
Set-CellText 117 4 "36,479 (94%)" "36,416 (94%)"
Set-CellText 117 5 "35,238 (94%)" "35,179 (94%)"
Set-CellText 117 6 "135,661 (95%)" "135,488 (95%)"

# CP1 row
Set-CellText 118 3 "49 (<0.1%)" "141 (<0.1%)"
Set-CellText 118 4 "8 (<0.1%)" "26 (<0.1%)"
Set-CellText 118 5 "10 (<0.1%)" "28 (<0.1%)"
Set-CellText 118 6 "31 (<0.1%)" "87 (<0.1%)"

# CP2 row
Set-CellText 119 3 "11,426 (5.2%)" "11,630 (5.3%)"
Set-CellText 119 4 "2,186 (5.7%)" "2,231 (5.8%)"
Set-CellText 119 5 "2,049 (5.5%)" "2,091 (5.6%)"
Set-CellText 119 6 "7,191 (5.0%)" "7,308 (5.1%)"

# Missing row (under CP block)
Set-CellText 120 3 "28,693" "28,692"
Set-CellText 120 5 "3,660" "3,659"

# t row (N)
Set-CellText 121 2 "218,853" "218,713"

# Median (IQR) row
Set-CellText 122 3 "236 (96, 416)" "236 (95, 416)"
Set-CellText 122 5 "455 (256, 591)" "455 (255, 590)"
Set-CellText 122 6 "210 (85, 368)" "209 (85, 367)"

# Missing row (under t block)
Set-CellText 124 3 "28,693" "28,833"
Set-CellText 124 4 "5,233" "5,259"
Set-CellText 124 5 "3,660" "3,687"
Set-CellText 124 6 "19,800" "19,887"

# incident_dm row (N + values)
Set-CellText 125 2 "218,853" "218,854"
Set-CellText 125 3 "11,475 (5.2%)" "11,771 (5.4%)"
Set-CellText 125 4 "2,194 (5.7%)" "2,257 (5.8%)"
Set-CellText 125 5 "2,059 (5.5%)" "2,119 (5.7%)"
Set-CellText 125 6 "7,222 (5.1%)" "7,395 (5.2%)"

# Missing row (under incident_dm block)
Set-CellText 126 3 "28,693" "28,692"
Set-CellText 126 5 "3,660" "3,659"

# SUPREMEDM Available row
Set-CellText 134 3 "218,853 (88%)" "218,854 (88%)"
Set-CellText 134 5 "37,297 (91%)" "37,298 (91%)"

# SUPREMEDM Unavailable row
Set-CellText 135 3 "28,693 (12%)" "28,692 (12%)"
Set-CellText 135 5 "3,660 (8.9%)" "3,659 (8.9%)"
